$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.455.92"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "1.955.96"
$ws.Range("E3").Value = "  -1.28%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'244.27"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").Value = "'0.611"
$ws.Range("E6").Value = "  -2.69%  "
$ws.Range("D7").Value = "'57.96"
$ws.Range("E7").Value = "  -3.69%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -3.50%  "
$ws.Range("D10").Value = "'0.0858"
$ws.Range("E10").Value = "  +5.18%  "
$ws.Range("E11").Value = "  +0.37%  "
$ws.Range("D12").Value = "2.242.23"
$ws.Range("E12").Value = "  -1.29%  "
$ws.Range("D13").Value = "'21.67"
$ws.Range("E13").Value = "  -6.47%  "
$ws.Range("E14").Value = "  -4.72%  "
$ws.Range("D15").Value = "'13.55"
$ws.Range("E15").Value = "  -4.75%  "
$ws.Range("E16").Value = "  -4.26%  "
$ws.Range("D17").Value = "1.956.57"
$ws.Range("E17").Value = "  -1.12%  "
$ws.Range("D18").Value = "36.399.12"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0886"
$ws.Range("E19").Value = "  +2.06%  "
$ws.Range("B20").Value = "Litecoin"
$ws.Range("C20").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D20").Value = "'70.05"
$ws.Range("E20").Value = "  -1.76%  "
$ws.Range("D21").Value = "'229.74"
$ws.Range("E21").Value = "  -2.30%  "
$ws.Range("D22").Value = "'5.09"
$ws.Range("E22").Value = "  -4.18%  "
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("E24").Value = "  -6.70%  "
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("D26").Value = "'9.31"
$ws.Range("E26").Value = "  -5.80%  "
$ws.Range("D27").Value = "'162.03"
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("D28").Value = "'0.136"
$ws.Range("E28").Value = "  +1.50%  "
$ws.Range("D29").Value = "'19.44"
$ws.Range("E29").Value = "  -2.23%  "
$ws.Range("D30").Value = "'0.119"
$ws.Range("E30").Value = "  -1.35%  "
$ws.Range("D31").Value = "'1.16"
$ws.Range("E31").Value = "  -0.90%  "
$ws.Range("D32").Value = "'4.67"
$ws.Range("E32").Value = "  -5.41%  "
$ws.Range("D33").Value = "'0.0648"
$ws.Range("E33").Value = "  +3.40%  "
$ws.Range("D34").Value = "'4.29"
$ws.Range("E34").Value = "  -4.49%  "
$ws.Range("D35").Value = "'6.31"
$ws.Range("E35").Value = "  +2.93%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").Value = "'1.79"
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("E38").Value = "  -5.28%  "
$ws.Range("D39").Value = "'3.09"
$ws.Range("E39").Value = "  -0.29%  "
$ws.Range("D40").Value = "'0.0984"
$ws.Range("E40").Value = "  +0.33%  "
$ws.Range("D41").Value = "'2.91"
$ws.Range("E41").Value = "  +0.73%  "
$ws.Range("E42").Value = "  -5.51%  "
$ws.Range("E43").Value = "  -1.29%  "
$ws.Range("D44").Value = "'15.81"
$ws.Range("E44").Value = "  -3.56%  "
$ws.Range("D45").Value = "1.353.54"
$ws.Range("E45").Value = "  -1.05%  "
$ws.Range("E46").Value = "  -6.08%  "
$ws.Range("D47").Value = "'88.05"
$ws.Range("E47").Value = "  -5.22%  "
$ws.Range("D48").Value = "'7.17"
$ws.Range("E48").Value = "  -6.20%  "
$ws.Range("D49").Value = "'2.84"
$ws.Range("E49").Value = "  -0.26%  "
$ws.Range("D50").Value = "'45.14"
$ws.Range("E50").Value = "  -2.25%  "
$ws.Range("D51").Value = "2.133.52"
$ws.Range("E51").Value = "  -1.43%  "
